# Append three new data rows (107-109) to the worksheet, matching the
# rows already present (unit/sex = THS_PERS / T) but with a geo value
# that is just a run of spaces (1, 2 and 3 spaces respectively) and
# years 2020, 2021, 2022. No values/rate_of_change are populated for
# these rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 107; Geo = " ";   Year = "2020" },
    @{ Row = 108; Geo = "  ";  Year = "2021" },
    @{ Row = 109; Geo = "   "; Year = "2022" }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    $ws.Cells.Item($r, 1).Value = $entry.Geo
    $ws.Cells.Item($r, 2).Value = "THS_PERS"
    $ws.Cells.Item($r, 3).Value = "T"

    # Force the year to be stored as text (matching the original file,
    # where it is an inline string, not a number) while keeping the
    # cell's style/number-format the same as a plain, unformatted cell.
    $ws.Cells.Item($r, 4).Value = "'" + $entry.Year
    $ws.Cells.Item($r, 4).Style = "Normal"
}
